# Case and Fatality Demographics Data Updated
# Apply updated case/fatality counts (one new data point) across the
# three "Fatalities by ..." sheets, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fatalities by Age")
$ws.Range("K8").Value = 408
$ws.Range("L8").Value = 1353
$ws.Range("K15").Value = 10778
$ws.Range("L15").Value = 31313
$ws.Range("K17").Value = 3576
$ws.Range("L17").Value = 9943
$ws.Range("G20").Value = 110
$ws.Range("L20").Value = 1245
$ws.Range("I23").Value = 275
$ws.Range("L23").Value = 1177
$ws.Range("H24").Value = 1297
$ws.Range("I24").Value = 1491
$ws.Range("K24").Value = 1271
$ws.Range("L24").Value = 6879
$ws.Range("E25").Value = 151
$ws.Range("H25").Value = 1533
$ws.Range("J25").Value = 1665
$ws.Range("L25").Value = 8293
$ws.Range("I26").Value = 1089
$ws.Range("J26").Value = 850
$ws.Range("L26").Value = 4087
$ws.Range("G27").Value = 162
$ws.Range("H27").Value = 311
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 356
$ws.Range("L27").Value = 1650
$ws.Range("E28").Value = 28
$ws.Range("F28").Value = 68
$ws.Range("G28").Value = 171
$ws.Range("H28").Value = 285
$ws.Range("I28").Value = 428
$ws.Range("J28").Value = 436
$ws.Range("K28").Value = 449
$ws.Range("L28").Value = 1871
$ws.Range("E29").Value = 551
$ws.Range("F29").Value = 1739
$ws.Range("G29").Value = 3865
$ws.Range("H29").Value = 6773
$ws.Range("I29").Value = 10318
$ws.Range("J29").Value = 10522
$ws.Range("K29").Value = 10494
$ws.Range("L29").Value = 44349
$ws.Range("E31").Value = 714
$ws.Range("F31").Value = 2251
$ws.Range("G31").Value = 5440
$ws.Range("H31").Value = 10289
$ws.Range("I31").Value = 16884
$ws.Range("J31").Value = 18688
$ws.Range("K31").Value = 21272
$ws.Range("L31").Value = 75662

$ws = $wb.Worksheets.Item("Fatalities by Sex")
$ws.Range("C8").Value = 847
$ws.Range("E8").Value = 1353
$ws.Range("C15").Value = 18121
$ws.Range("B17").Value = 4186
$ws.Range("E17").Value = 9943
$ws.Range("C20").Value = 693
$ws.Range("E20").Value = 1245
$ws.Range("B23").Value = 495
$ws.Range("E23").Value = 1177
$ws.Range("B24").Value = 2804
$ws.Range("C24").Value = 4075
$ws.Range("E24").Value = 6879
$ws.Range("B25").Value = 3438
$ws.Range("C25").Value = 4855
$ws.Range("E25").Value = 8293
$ws.Range("B26").Value = 1788
$ws.Range("C26").Value = 2299
$ws.Range("E26").Value = 4087
$ws.Range("C27").Value = 965
$ws.Range("E27").Value = 1650
$ws.Range("B28").Value = 782
$ws.Range("C28").Value = 1089
$ws.Range("E28").Value = 1871
$ws.Range("B29").Value = 18502
$ws.Range("C29").Value = 25847
$ws.Range("E29").Value = 44349
$ws.Range("B31").Value = 31693
$ws.Range("C31").Value = 43968
$ws.Range("E31").Value = 75661

$ws = $wb.Worksheets.Item("Fatalities by RaceEth")
$ws.Range("F8").Value = 351
$ws.Range("H8").Value = 1353
$ws.Range("F15").Value = 11892
$ws.Range("H15").Value = 31313
$ws.Range("F17").Value = 4962
$ws.Range("H17").Value = 9943
$ws.Range("F20").Value = 450
$ws.Range("H20").Value = 1245
$ws.Range("C23").Value = 132
$ws.Range("H23").Value = 1177
$ws.Range("D24").Value = 2541
$ws.Range("F24").Value = 3402
$ws.Range("H24").Value = 6879
$ws.Range("B25").Value = 98
$ws.Range("C25").Value = 1004
$ws.Range("D25").Value = 2979
$ws.Range("F25").Value = 4144
$ws.Range("H25").Value = 8293
$ws.Range("F26").Value = 1976
$ws.Range("H26").Value = 4087
$ws.Range("C27").Value = 141
$ws.Range("D27").Value = 669
$ws.Range("F27").Value = 804
$ws.Range("H27").Value = 1650
$ws.Range("C28").Value = 153
$ws.Range("D28").Value = 721
$ws.Range("F28").Value = 948
$ws.Range("H28").Value = 1871
$ws.Range("B29").Value = 807
$ws.Range("C29").Value = 4854
$ws.Range("D29").Value = 17403
$ws.Range("F29").Value = 20964
$ws.Range("H29").Value = 44349
$ws.Range("B31").Value = 1384
$ws.Range("C31").Value = 8031
$ws.Range("D31").Value = 32881
$ws.Range("F31").Value = 32856
$ws.Range("H31").Value = 75662

# Cosmetic: "Fatalities by Age" tab zoom normalized to 100% in the saved view
$ws1 = $wb.Worksheets.Item("Fatalities by Age")
$ws1.Activate()
$excel.ActiveWindow.Zoom = 100
